# Collapse. expand sidebar menu, filter popup
#
# The underlying data table had no header row and included two rows
# ("3" / sdvdsvsx / ... / Ngung kinh doanh, and "5" / fhfdgdf / ... /
# Ngung kinh doanh) that the filter popup should no longer surface.
# This adds a header row above the data and removes those two rows,
# leaving products 1, 2 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; everything else shifts down one row.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Mã Sp"
$ws.Range("B1").Value = "Tên Sp"
$ws.Range("C1").Value = "Số lượng"
$ws.Range("D1").Value = "Mô tả"
$ws.Range("E1").Value = "Tình trạng"

# After the insert, the original 5 data rows now live in rows 2-6:
#   row 2 -> "1"
#   row 3 -> "2"
#   row 4 -> "3"   (sdvdsvsx / 678 / dfvfv / Ngung kinh doanh)  -- remove
#   row 5 -> "4"
#   row 6 -> "5"   (fhfdgdf / 0 / sdfsd / Ngung kinh doanh)      -- remove
#
# Delete row 4 first; row 6 then becomes row 5.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(5).Delete()
